# testing/performance/perfData.xlsx
# pidan: test to get datas to write summary
#
# Add a new row (row 3) of benchmark data and update the "tune mc nc kc"
# result in row 2 (column G) with a newer measurement.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values are written in the same order the author's workbook records
# them in the shared string table: G3, then C3, then G2.
$ws.Range("A3").Value = 1152
$ws.Range("G3").Value = "67(88%)"
$ws.Range("C3").Value = "66(86%)"
$ws.Range("G2").Value = "65(85%)"

# Leave the selection on the last cell that was edited.
$ws.Range("G3").Select() | Out-Null
